$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff marks several additional task checkboxes as completed
# (checkmark "ü" in Wingdings, with the matching "checked" cell style
# taken from the fill/border pattern already used elsewhere in the
# same column). Row 4 already carries the "checked" formatting for
# each column, so we copy its per-column format onto the target cells
# before writing the checkmark value.

$targets = @(
    @{ Cell = "D9";  FormatSource = "D4" },
    @{ Cell = "E9";  FormatSource = "E4" },
    @{ Cell = "F9";  FormatSource = "F4" },
    @{ Cell = "I23"; FormatSource = "I4" },
    @{ Cell = "H25"; FormatSource = "H4" },
    @{ Cell = "G28"; FormatSource = "G4" }
)

foreach ($t in $targets) {
    $ws.Range($t.FormatSource).Copy()
    $ws.Range($t.Cell).PasteSpecial(-4122)
    $ws.Range($t.Cell).Value = "ü"
}

$excel.CutCopyMode = 0
